$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'245.77"
$ws.Range("D2").Style = "Normal"
$ws.Range("F2").Value = "'31-12-2022"
$ws.Range("F2").Style = "Normal"
$ws.Range("G2").Value = "'0"
$ws.Range("G2").Style = "Normal"
$ws.Range("D3").Value = "'25.44"
$ws.Range("D3").Style = "Normal"
$ws.Range("F3").Value = "'31-12-2022"
$ws.Range("F3").Style = "Normal"
$ws.Range("G3").Value = "'0"
$ws.Range("G3").Style = "Normal"
$ws.Range("D4").Value = "'5.033"
$ws.Range("D4").Style = "Normal"
$ws.Range("F4").Value = "'31-12-2022"
$ws.Range("F4").Style = "Normal"
$ws.Range("G4").Value = "'0"
$ws.Range("G4").Style = "Normal"
$ws.Range("F5").Value = "'31-12-2022"
$ws.Range("F5").Style = "Normal"
$ws.Range("G5").Value = "'0"
$ws.Range("G5").Style = "Normal"
$ws.Range("D6").Value = "'6.561"
$ws.Range("D6").Style = "Normal"
$ws.Range("F6").Value = "'31-12-2022"
$ws.Range("F6").Style = "Normal"
$ws.Range("G6").Value = "'0"
$ws.Range("G6").Style = "Normal"
$ws.Range("D7").Value = "'3.021"
$ws.Range("D7").Style = "Normal"
$ws.Range("F7").Value = "'31-12-2022"
$ws.Range("F7").Style = "Normal"
$ws.Range("G7").Value = "'0"
$ws.Range("G7").Style = "Normal"
$ws.Range("D8").Value = "'0.8174"
$ws.Range("D8").Style = "Normal"
$ws.Range("F8").Value = "'31-12-2022"
$ws.Range("F8").Style = "Normal"
$ws.Range("G8").Value = "'0"
$ws.Range("G8").Style = "Normal"
$ws.Range("D9").Value = "'0.8363"
$ws.Range("D9").Style = "Normal"
$ws.Range("F9").Value = "'31-12-2022"
$ws.Range("F9").Style = "Normal"
$ws.Range("G9").Value = "'0"
$ws.Range("G9").Style = "Normal"
$ws.Range("D10").Value = "'0.0005950"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "9OneONE"
$ws.Range("F10").Value = "'31-12-2022"
$ws.Range("F10").Style = "Normal"
$ws.Range("G10").Value = "'0"
$ws.Range("G10").Style = "Normal"
$ws.Range("D11").Value = "'0.1341"
$ws.Range("D11").Style = "Normal"
$ws.Range("F11").Value = "'31-12-2022"
$ws.Range("F11").Style = "Normal"
$ws.Range("G11").Value = "'0"
$ws.Range("G11").Style = "Normal"
$ws.Range("D12").Value = "'0.06957"
$ws.Range("D12").Style = "Normal"
$ws.Range("F12").Value = "'31-12-2022"
$ws.Range("F12").Style = "Normal"
$ws.Range("G12").Value = "'0"
$ws.Range("G12").Style = "Normal"
$ws.Range("B13").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C13").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D13").Value = "'0.03215"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "12LiechtensteinCryptoassetsExchangeLCX"
$ws.Range("F13").Value = "'31-12-2022"
$ws.Range("F13").Style = "Normal"
$ws.Range("G13").Value = "'0"
$ws.Range("G13").Style = "Normal"
$ws.Range("B14").Value = "BitrueCoin"
$ws.Range("C14").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D14").Value = "'0.02820"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "13BitrueCoinBTR"
$ws.Range("F14").Value = "'31-12-2022"
$ws.Range("F14").Style = "Normal"
$ws.Range("G14").Value = "'0"
$ws.Range("G14").Style = "Normal"
$ws.Range("B15").Value = "BitMartToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D15").Value = "'0.09396"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "14BitMartTokenBMX"
$ws.Range("F15").Value = "'31-12-2022"
$ws.Range("F15").Style = "Normal"
$ws.Range("G15").Value = "'0"
$ws.Range("G15").Style = "Normal"
$ws.Range("B16").Value = "BitForexToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D16").Value = "'0.001516"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "15BitForexTokenBF"
$ws.Range("F16").Value = "'31-12-2022"
$ws.Range("F16").Style = "Normal"
$ws.Range("G16").Value = "'0"
$ws.Range("G16").Style = "Normal"
$ws.Range("B17").Value = "TigerCash"
$ws.Range("C17").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D17").Value = "'0.006147"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "16TigerCashTCH"
$ws.Range("F17").Value = "'31-12-2022"
$ws.Range("F17").Style = "Normal"
$ws.Range("G17").Value = "'0"
$ws.Range("G17").Style = "Normal"
$ws.Range("B18").Value = "LEO"
$ws.Range("C18").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D18").Value = "'3.507"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "17LEOLEO"
$ws.Range("F18").Value = "'31-12-2022"
$ws.Range("F18").Style = "Normal"
$ws.Range("G18").Value = "'0"
$ws.Range("G18").Style = "Normal"
$ws.Range("B19").Value = "BTSEToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D19").Value = "'2.091"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "18BTSETokenBTSE"
$ws.Range("F19").Value = "'31-12-2022"
$ws.Range("F19").Style = "Normal"
$ws.Range("G19").Value = "'0"
$ws.Range("G19").Style = "Normal"
$ws.Range("B20").Value = "BitpandaEcosystemToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D20").Value = "'0.3179"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "19BitpandaEcosystemTokenBEST"
$ws.Range("F20").Value = "'31-12-2022"
$ws.Range("F20").Style = "Normal"
$ws.Range("G20").Value = "'0"
$ws.Range("G20").Style = "Normal"
$ws.Range("D21").Value = "'0.1320"
$ws.Range("D21").Style = "Normal"
$ws.Range("F21").Value = "'31-12-2022"
$ws.Range("F21").Style = "Normal"
$ws.Range("G21").Value = "'0"
$ws.Range("G21").Style = "Normal"
$ws.Range("D22").Value = "'3.739"
$ws.Range("D22").Style = "Normal"
$ws.Range("F22").Value = "'31-12-2022"
$ws.Range("F22").Style = "Normal"
$ws.Range("G22").Value = "'0"
$ws.Range("G22").Style = "Normal"
$ws.Range("D23").Value = "'0.04690"
$ws.Range("D23").Style = "Normal"
$ws.Range("F23").Value = "'31-12-2022"
$ws.Range("F23").Style = "Normal"
$ws.Range("G23").Value = "'0"
$ws.Range("G23").Style = "Normal"
$ws.Range("F24").Value = "'31-12-2022"
$ws.Range("F24").Style = "Normal"
$ws.Range("G24").Value = "'0"
$ws.Range("G24").Style = "Normal"
$ws.Range("D25").Value = "'0.001242"
$ws.Range("D25").Style = "Normal"
$ws.Range("F25").Value = "'31-12-2022"
$ws.Range("F25").Style = "Normal"
$ws.Range("G25").Value = "'0"
$ws.Range("G25").Style = "Normal"
$ws.Range("D26").Value = "'0.004285"
$ws.Range("D26").Style = "Normal"
$ws.Range("F26").Value = "'31-12-2022"
$ws.Range("F26").Style = "Normal"
$ws.Range("G26").Value = "'0"
$ws.Range("G26").Style = "Normal"
$ws.Range("F27").Value = "'31-12-2022"
$ws.Range("F27").Style = "Normal"
$ws.Range("G27").Value = "'0"
$ws.Range("G27").Style = "Normal"
$ws.Range("D28").Value = "'0.0001385"
$ws.Range("D28").Style = "Normal"
$ws.Range("F28").Value = "'31-12-2022"
$ws.Range("F28").Style = "Normal"
$ws.Range("G28").Value = "'0"
$ws.Range("G28").Style = "Normal"
$ws.Range("F29").Value = "'31-12-2022"
$ws.Range("F29").Style = "Normal"
$ws.Range("G29").Value = "'0"
$ws.Range("G29").Style = "Normal"
$ws.Range("F30").Value = "'31-12-2022"
$ws.Range("F30").Style = "Normal"
$ws.Range("G30").Value = "'0"
$ws.Range("G30").Style = "Normal"
$ws.Range("F31").Value = "'31-12-2022"
$ws.Range("F31").Style = "Normal"
$ws.Range("G31").Value = "'0"
$ws.Range("G31").Style = "Normal"
$ws.Range("F32").Value = "'31-12-2022"
$ws.Range("F32").Style = "Normal"
$ws.Range("G32").Value = "'0"
$ws.Range("G32").Style = "Normal"
$ws.Range("F33").Value = "'31-12-2022"
$ws.Range("F33").Style = "Normal"
$ws.Range("G33").Value = "'0"
$ws.Range("G33").Style = "Normal"
$ws.Range("F34").Value = "'31-12-2022"
$ws.Range("F34").Style = "Normal"
$ws.Range("G34").Value = "'0"
$ws.Range("G34").Style = "Normal"
$ws.Range("F35").Value = "'31-12-2022"
$ws.Range("F35").Style = "Normal"
$ws.Range("G35").Value = "'0"
$ws.Range("G35").Style = "Normal"
$ws.Range("F36").Value = "'31-12-2022"
$ws.Range("F36").Style = "Normal"
$ws.Range("G36").Value = "'0"
$ws.Range("G36").Style = "Normal"
$ws.Range("F37").Value = "'31-12-2022"
$ws.Range("F37").Style = "Normal"
$ws.Range("G37").Value = "'0"
$ws.Range("G37").Style = "Normal"
$ws.Range("F38").Value = "'31-12-2022"
$ws.Range("F38").Style = "Normal"
$ws.Range("G38").Value = "'0"
$ws.Range("G38").Style = "Normal"
$ws.Range("F39").Value = "'31-12-2022"
$ws.Range("F39").Style = "Normal"
$ws.Range("G39").Value = "'0"
$ws.Range("G39").Style = "Normal"
$ws.Range("D40").Value = "'0.03654"
$ws.Range("D40").Style = "Normal"
$ws.Range("F40").Value = "'31-12-2022"
$ws.Range("F40").Style = "Normal"
$ws.Range("G40").Value = "'0"
$ws.Range("G40").Style = "Normal"
$ws.Range("D41").Value = "'0.006180"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "40KickTokenKICKBestin24h"
$ws.Range("F41").Value = "'31-12-2022"
$ws.Range("F41").Style = "Normal"
$ws.Range("G41").Value = "'0"
$ws.Range("G41").Style = "Normal"
$ws.Range("F42").Value = "'31-12-2022"
$ws.Range("F42").Style = "Normal"
$ws.Range("G42").Value = "'0"
$ws.Range("G42").Style = "Normal"
$ws.Range("F43").Value = "'31-12-2022"
$ws.Range("F43").Style = "Normal"
$ws.Range("G43").Value = "'0"
$ws.Range("G43").Style = "Normal"
$ws.Range("D44").Value = "'0.008211"
$ws.Range("D44").Style = "Normal"
$ws.Range("F44").Value = "'31-12-2022"
$ws.Range("F44").Style = "Normal"
$ws.Range("G44").Value = "'0"
$ws.Range("G44").Style = "Normal"
$ws.Range("D45").Value = "'0.00005297"
$ws.Range("D45").Style = "Normal"
$ws.Range("F45").Value = "'31-12-2022"
$ws.Range("F45").Style = "Normal"
$ws.Range("G45").Value = "'0"
$ws.Range("G45").Style = "Normal"
$ws.Range("F46").Value = "'31-12-2022"
$ws.Range("F46").Style = "Normal"
$ws.Range("G46").Value = "'0"
$ws.Range("G46").Style = "Normal"
$ws.Range("F47").Value = "'31-12-2022"
$ws.Range("F47").Style = "Normal"
$ws.Range("G47").Value = "'0"
$ws.Range("G47").Style = "Normal"
$ws.Range("F48").Value = "'31-12-2022"
$ws.Range("F48").Style = "Normal"
$ws.Range("G48").Value = "'0"
$ws.Range("G48").Style = "Normal"
$ws.Range("F49").Value = "'31-12-2022"
$ws.Range("F49").Style = "Normal"
$ws.Range("G49").Value = "'0"
$ws.Range("G49").Style = "Normal"
$ws.Range("F50").Value = "'31-12-2022"
$ws.Range("F50").Style = "Normal"
$ws.Range("G50").Value = "'0"
$ws.Range("G50").Style = "Normal"
$ws.Range("F51").Value = "'31-12-2022"
$ws.Range("F51").Style = "Normal"
$ws.Range("G51").Value = "'0"
$ws.Range("G51").Style = "Normal"
